$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.732.41'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.856.94'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9988'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.73'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6412'
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9994'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.31'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +4.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07522'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2983'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '24.48'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +5.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07667'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.867.25'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.038'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6920'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.87'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009834'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +9.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.089'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '29.743.12'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.116.94'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '236.35'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.68'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9996'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.536'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.9991'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.90'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1422'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.535'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.92'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06209'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +6.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.492'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.282'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.164'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.103'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.895'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.62%  '
$ws.Range("E36").Value = '  +2.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7290'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.603'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.827'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01787'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.204.14'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9246'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.230'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.030.27'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.18%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9996'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.99'
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '66.62'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000119'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.06%  '
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4061'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.93%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.187'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05793'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.88%  '
